# Apply cryptos list update (price & volume refresh) — commit: Updated cryptos list on Tue Mar 21 15:34:32 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the cell to keep an exact text representation (e.g. "45.60", "1.003")
    # instead of Excel's automatic numeric coercion, then restore the default style
    # so the cell XML matches the original (un-styled) inline-string cells.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# --- Regular per-row Price (D) / Volume(1h) (E) updates ---
Set-TextCell $ws.Range("D2") "28.038.26"
$ws.Range("E2").Value = "  -0.30%  "
Set-TextCell $ws.Range("D3") "1.815.86"
$ws.Range("E3").Value = "  +2.13%  "
Set-TextCell $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextCell $ws.Range("D5") "337.44"
$ws.Range("E5").Value = "  -0.52%  "
Set-TextCell $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  -0.01%  "
Set-TextCell $ws.Range("D7") "0.4268"
$ws.Range("E7").Value = "  +11.60%  "
Set-TextCell $ws.Range("D8") "0.3505"
$ws.Range("E8").Value = "  +2.63%  "
Set-TextCell $ws.Range("D9") "45.60"
$ws.Range("E9").Value = "  -2.46%  "
Set-TextCell $ws.Range("D10") "1.148"
$ws.Range("E10").Value = "  +0.58%  "
Set-TextCell $ws.Range("D11") "0.07451"
$ws.Range("E11").Value = "  +0.90%  "
Set-TextCell $ws.Range("D12") "22.95"
$ws.Range("E12").Value = "  -1.36%  "
Set-TextCell $ws.Range("D13") "1.002"
$ws.Range("E13").Value = "  +0.19%  "
Set-TextCell $ws.Range("D14") "6.262"
$ws.Range("E14").Value = "  -1.84%  "
Set-TextCell $ws.Range("D15") "7.290"
$ws.Range("E15").Value = "  -1.37%  "
Set-TextCell $ws.Range("D16") "1.811.32"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("E17").Value = "  +0.93%  "
Set-TextCell $ws.Range("D18") "0.06689"
$ws.Range("E18").Value = "  +0.46%  "
Set-TextCell $ws.Range("D19") "81.96"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("E20").Value = "  +0.04%  "
Set-TextCell $ws.Range("D23") "28.057.40"
$ws.Range("E23").Value = "  -0.21%  "
Set-TextCell $ws.Range("D24") "12.06"
$ws.Range("E24").Value = "  +0.08%  "
Set-TextCell $ws.Range("D25") "2.390"
$ws.Range("E25").Value = "  +0.41%  "
Set-TextCell $ws.Range("D26") "2.486"
$ws.Range("E26").Value = "  +3.13%  "
Set-TextCell $ws.Range("D27") "20.73"
$ws.Range("E27").Value = "  +0.13%  "
Set-TextCell $ws.Range("D28") "156.11"
$ws.Range("E28").Value = "  +1.21%  "
Set-TextCell $ws.Range("D29") "2.018.65"
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("E30").Value = "  -10.46%  "
Set-TextCell $ws.Range("D31") "132.47"
$ws.Range("E31").Value = "  -1.40%  "
Set-TextCell $ws.Range("D32") "4.064"
$ws.Range("E32").Value = "  +1.10%  "
Set-TextCell $ws.Range("D33") "5.951"
$ws.Range("E33").Value = "  -1.77%  "
Set-TextCell $ws.Range("D34") "0.09213"
$ws.Range("E34").Value = "  +3.41%  "
Set-TextCell $ws.Range("D35") "12.36"
$ws.Range("E35").Value = "  -2.65%  "
Set-TextCell $ws.Range("D36") "0.02357"
$ws.Range("E36").Value = "  -2.09%  "
Set-TextCell $ws.Range("D37") "0.6724"
$ws.Range("E37").Value = "  -1.50%  "
Set-TextCell $ws.Range("D38") "5.244"
$ws.Range("E38").Value = "  -0.83%  "
Set-TextCell $ws.Range("D39") "0.06273"
$ws.Range("E39").Value = "  -1.61%  "
Set-TextCell $ws.Range("D40") "0.2172"
$ws.Range("E40").Value = "  +0.63%  "
Set-TextCell $ws.Range("D41") "1.496"
$ws.Range("E41").Value = "  +0.05%  "
Set-TextCell $ws.Range("D42") "1.217"
$ws.Range("E42").Value = "  -1.62%  "
Set-TextCell $ws.Range("D43") "8.101"
$ws.Range("E43").Value = "  -1.29%  "
Set-TextCell $ws.Range("D46") "3.873"
$ws.Range("E46").Value = "  +0.26%  "
Set-TextCell $ws.Range("D47") "0.6125"
$ws.Range("E47").Value = "  -2.19%  "
Set-TextCell $ws.Range("D48") "128.37"
$ws.Range("E48").Value = "  -3.71%  "
Set-TextCell $ws.Range("D49") "2.046"
$ws.Range("E49").Value = "  -1.04%  "
Set-TextCell $ws.Range("D50") "1.180"
$ws.Range("E50").Value = "  -2.46%  "
Set-TextCell $ws.Range("D51") "0.07110"
$ws.Range("E51").Value = "  -5.22%  "

# --- Rows 21/22 swapped places: Avalanche <-> Uniswap (new prices/volumes) ---
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell $ws.Range("D21") "6.457"
$ws.Range("E21").Value = "  +0.89%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell $ws.Range("D22") "17.25"
$ws.Range("E22").Value = "  -0.69%  "

# --- Rows 44/45 swapped places: EnergySwap <-> Frax (new prices/volumes) ---
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell $ws.Range("D44") "0.9998"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D45") "13.99"
$ws.Range("E45").Value = "  -1.62%  "
